$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted ahead of the existing
# "Zapallo italiano" series (row 319), shifting every following row down
# by one and growing the used range to A1:R396.
$ws.Rows.Item(319).Insert()

$ws.Cells.Item(319, 1).Value = 5
$ws.Cells.Item(319, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(319, 3).Value = "Maule"
$ws.Cells.Item(319, 4).Value = 44782
$ws.Cells.Item(319, 5).Value = 7
$ws.Cells.Item(319, 6).Value = 100112032
$ws.Cells.Item(319, 7).Value = "Zapallo italiano"
$ws.Cells.Item(319, 8).Value = "Sin especificar"
$ws.Cells.Item(319, 9).Value = "Primera"
$ws.Cells.Item(319, 10).Value = 300
$ws.Cells.Item(319, 11).Value = 20000
$ws.Cells.Item(319, 12).Value = 20000
$ws.Cells.Item(319, 13).Value = 20000
$ws.Cells.Item(319, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(319, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(319, 16).Value = 400
$ws.Cells.Item(319, 17).Value = 50
$ws.Cells.Item(319, 18).Value = "Hortaliza"
